# Update the "drop in" user rows on the home sheet so the UI can tell
# whether a visitor needs to login/signup or is already dropped into a game.
#
# Row 2 (Nic Bolton) keeps the same person but the games/hosted-games
# bookkeeping columns are refreshed, and the password is normalized to a
# text value.
# Row 3's former placeholder user ("Steve Jobs") is replaced with a real
# second account (Albert Bolton) including a new phone/password and
# matching games state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Nic Bolton -------------------------------------------------
# password needs to be stored as the text "1234" (not the number 1234).
# Build it as a text value via TEXT(), then copy/paste-special as values so
# the destination cell keeps its existing (default) formatting.
$ws.Range("Z1").Formula = "=TEXT(1234,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# games / hosted games columns
$ws.Range("F2").Value = "1;"
$ws.Range("G2").Value = "0;1;"

# --- Row 3: Albert Bolton (was Steve Jobs) -----------------------------
$ws.Range("B3").Value = "Albert Bolton"
$ws.Range("C3").Value = "albabolton@me.com"
$ws.Range("D3").Value = 5195663730
$ws.Range("E3").Value = "Ryder0498`$"
$ws.Range("F3").Value = "1;"
$ws.Range("G3").Value = ""
